# Apply "more mutant analysis with botium (TBC)" changes to results_pizza.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header in column G: "#tests = 52" ---
$ws.Cells.Item(1, 7).Value = "#tests = 52"
$ws.Cells.Item(1, 7).Font.Bold = $true
$ws.Cells.Item(1, 7).Interior.Color = 65535
$ws.Cells.Item(1, 7).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(1, 7).Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Cells.Item(1, 7).Borders.Item(10).LineStyle = 1  # xlEdgeRight

# Widen column F/G like the new col width for column F (65.71)
$ws.Columns.Item(6).ColumnWidth = 65.7109375

# --- Fill in newly-tracked mutant rows (26, 28-32) ---
$ws.Cells.Item(26, 2).Value = "Done"
$ws.Cells.Item(26, 5).Value = "survived"

$ws.Cells.Item(28, 2).Value = "Done"
$ws.Cells.Item(28, 5).Value = "KILLED"
$ws.Cells.Item(28, 6).Value = "Hace que salte un intent distinto (expected toppings but found order drinks)"

$ws.Cells.Item(29, 2).Value = "Done"
$ws.Cells.Item(29, 5).Value = "KILLED"
$ws.Cells.Item(29, 6).Value = "Hace que salte un intent distinto"

$ws.Cells.Item(30, 2).Value = "Done"
$ws.Cells.Item(30, 5).Value = "survived"

$ws.Cells.Item(31, 2).Value = "Done"
$ws.Cells.Item(31, 5).Value = "survived"

$ws.Cells.Item(32, 2).Value = "Done"
$ws.Cells.Item(32, 5).Value = "survived"

# --- Mut Score label + formula in H4/I4 ---
$ws.Cells.Item(4, 8).Value = "Mut Score"
$ws.Cells.Item(4, 9).Formula = '=COUNTIF(E2:E94, "Killed")/(COUNTIF(E2:E94, "Killed")+COUNTIF(E2:E94, "survived"))'

# --- Update selection to match the saved view state ---
$ws.Range("F32").Select()
